{"js": "// Apply the two localized text edits described by the diff:\n//   1. \"q<page>\"            -> \"<page>\"        (drop the stray leading \"q\")\n//   2. \"moing\"               -> \"moingdre\"       (append \"dre\")\n//   3. \"<ill/></ms>\"         -> \"</ms>\"          (drop the \"<ill/>\" tag)\n//\n// Each run in this document holds a short literal snippet of text, so a\n// scoped, case-sensitive search for each exact snippet identifies the\n// correct run; insertText(..., \"Replace\") then rewrites just that run's\n// text while preserving its run formatting (rPr).\n\nconst body = context.document.body;\n\nconst pageHits = body.search(\"q<page>\", { matchCase: true, matchWholeWord: false });\npageHits.load(\"items\");\n\nconst moingHits = body.search(\"moing\", { matchCase: true, matchWholeWord: false });\nmoingHits.load(\"items\");\n\nconst illHits = body.search(\"<ill/></ms>\", { matchCase: true, matchWholeWord: false });\nillHits.load(\"items\");\n\nawait context.sync();\n\nif (pageHits.items.length > 0) {\n  pageHits.items[0].insertText(\"<page>\", \"Replace\");\n}\n\nif (moingHits.items.length > 0) {\n  moingHits.items[0].insertText(\"moingdre\", \"Replace\");\n}\n\nif (illHits.items.length > 0) {\n  illHits.items[0].insertText(\"</ms>\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the three localized text edits described by the diff:\n#   1. \"q<page>\"      -> \"<page>\"    (drop the stray leading \"q\")\n#   2. \"moing\"        -> \"moingdre\"  (append \"dre\")\n#   3. \"<ill/></ms>\"  -> \"</ms>\"     (drop the \"<ill/>\" tag)\n#\n# Each snippet is a unique, literal run of text in the document, so a\n# plain (non-wildcard) Find/Replace targeted at Document.Content finds the\n# exact run and rewrites just its text, leaving every other run untouched.\n\n$wdFindContinue = 1\n$wdReplaceOne   = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-OneMatch([string]$findText, [string]$replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne) | Out-Null\n}\n\nReplace-OneMatch \"q<page>\" \"<page>\"\nReplace-OneMatch \"moing\" \"moingdre\"\nReplace-OneMatch \"<ill/></ms>\" \"</ms>\"\n\nWrite-Output \"done\"\n"}
